$d = $word.ActiveDocument

# 1. Eko-CORE bullet: remove " device", change "saving cost by half for patients with
#    fistula (AVF)" to "saving monthly cost by $900 for patients on dialysis"
$d.Content.Find.Execute("stethoscope device, saving cost by half for patients with fistula (AVF)", $true, $false, $false, $false, $false, $true, 1, $false, "stethoscope, saving monthly cost by `$900 for patients on dialysis", 2)

# 2. stenosis / AV fistula bullet
$d.Content.Find.Execute("to detect stenosis caused by AV fistula", $true, $false, $false, $false, $false, $true, 1, $false, "to detect vascular stenosis in arteriovenous fistula", 2)

# 3. machine learning models bullet
$d.Content.Find.Execute("(acc: 73.68%, AUC: 0.85) detecting stenosis", $true, $false, $false, $false, $false, $true, 1, $false, "(acc: 73.68%, AUC: 0.85) for stenosis detection", 2)

# 4. SAS / EDA bullet
$d.Content.Find.Execute("and performed EDA using dplyr", $true, $false, $false, $false, $false, $true, 1, $false, "and manipulated data using dplyr", 2)

# 5. GPA lines
$d.Content.Find.Execute("GPA: 3.6/4.0", $true, $false, $false, $false, $false, $true, 1, $false, "GPA: 3.6", 2)
$d.Content.Find.Execute("GPA: 3.8/4.0", $true, $false, $false, $false, $false, $true, 1, $false, "GPA: 3.8", 2)
